$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 271, pushing existing rows 271:297 down to 272:298.
$ws.Rows(271).Insert()

# Populate the newly-inserted row 271 with the new data point.
$ws.Range("A271").Value = 3
$ws.Range("B271").Value = "Femacal de La Calera"
$ws.Range("C271").Value = "Coquimbo"
$ws.Range("D271").Value = 45166
$ws.Range("E271").Value = 5
$ws.Range("F271").Value = 100112026
$ws.Range("G271").Value = "Haba"
$ws.Range("H271").Value = "Sin especificar"
$ws.Range("I271").Value = "Primera"
$ws.Range("J271").Value = 50
$ws.Range("K271").Value = 15000
$ws.Range("L271").Value = 15000
$ws.Range("M271").Value = 15000
$ws.Range("N271").Value = "$/saco 25 kilos"
$ws.Range("O271").Value = "Provincia de Limarí"
$ws.Range("P271").Value = 600
$ws.Range("Q271").Value = 25
$ws.Range("R271").Value = "Hortaliza"
